# Update "Login Credentials" workbook:
#  - Sheet1 row 3 & 4: login email changes from clarinenyauncho@gmail.com
#    to stella.ireri@tezzasolutions.com (value + hyperlink target)
#  - Sheet1 row 4: password changes from "changeme" to "!qwerty123"
#  - Sheet1 becomes the active sheet/tab (was Sheet2)
#  - Selection bookmarks updated: Sheet1 -> C18, Sheet2 -> A2

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- update the login email on rows 3 and 4 ---
$ws1.Range("B3").Value = "stella.ireri@tezzasolutions.com"
$ws1.Range("B4").Value = "stella.ireri@tezzasolutions.com"

# --- update the password on row 4 ---
$ws1.Range("C4").Value = "!qwerty123"

# --- repoint the mailto hyperlinks for the cells whose email changed ---
foreach ($h in $ws1.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$B$3' -or $addr -eq '$B$4') {
        $h.Address = "mailto:stella.ireri@tezzasolutions.com"
    }
}

# --- restore per-sheet selection bookmarks ---
# Set Sheet2's remembered selection first (it will not remain the active sheet)...
[void]$ws2.Range("A2").Select()
# ...then select on Sheet1 last so Sheet1 ends up the active/tabSelected sheet.
[void]$ws1.Range("C18").Select()
[void]$ws1.Activate()
